# Update "想去人数" (F column) figures for both the "展览" sheet and the
# "全部类型" sheet, which contain duplicated data for the same events.

$wb = $excel.ActiveWorkbook

# Map: row number -> new F value, for sheet "展览"
$sheet1Updates = @{
    2  = 77
    3  = 54
    4  = 4626
    5  = 1838
    8  = 3107
    11 = 261
    12 = 625
    13 = 534
    14 = 525
    15 = 375
    16 = 133
    18 = 1327
    19 = 123
    20 = 1596
    23 = 4
    25 = 531
    30 = 21
    32 = 3795
    33 = 758
    34 = 67
    35 = 686
    36 = 57
    37 = 1809
}

# Map: row number -> new F value, for sheet "全部类型"
$sheet4Updates = @{
    2  = 77
    3  = 54
    4  = 4626
    5  = 1838
    8  = 3107
    11 = 261
    12 = 625
    13 = 534
    14 = 525
    16 = 375
    17 = 133
    19 = 1328
    20 = 123
    21 = 1596
    24 = 4
    26 = 531
    31 = 21
    33 = 3795
    35 = 758
    36 = 67
    37 = 686
    38 = 57
    39 = 1809
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
